$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new label in column A, row 6
$ws.Range("A6").Value = "SNR"

# Add SNR formula for B6 (single, standalone formula)
$ws.Range("B6").Formula = "= 10 *LOG10(((B5-B3)^2)/((B4)^2))"

# Add SNR formula for C6:G6 as a shared formula, same as Excel's fill-right gesture
$ws.Range("C6:G6").Formula = "= 10 *LOG10(((C5-C3)^2)/((C4)^2))"

# Update selection to match the diff
$ws.Range("P26").Select()
